# Update the "Periodo Mora" column (E16:E19) on the account-statement
# worksheet so the periods read in ascending order (2303, 2304, 2305, 2306)
# instead of the previous descending order (2306, 2305, 2304, 2303).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2303"
$ws.Range("E17").Value = "2304"
$ws.Range("E18").Value = "2305"
$ws.Range("E19").Value = "2306"
